$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bio/common-name/verification data (rows 2-11) ---
$ws.Range("E2").Value = "N"

$ws.Range("B3").Value = "Kuranda Tree Frog"
$ws.Range("E3").Value = "Y"
$ws.Range("H3").Value = "Common Name?"

$ws.Range("B4").Value = "Peron's Tree Frog"
$ws.Range("H4").Value = "Common Name?"

$ws.Range("H5").Value = "Y"

$ws.Range("H6").Value = "Y"

$ws.Range("B7").Value = "Orange Thighed Tree Frog"
$ws.Range("H7").Value = "Common Name?"

$ws.Range("H8").Value = "Y"

$ws.Range("H9").Value = "Y"

$ws.Range("H10").Value = "Y"

$ws.Range("B11").Value = "Purple-crowned fairy wren"
$ws.Range("H11").Value = "Common Name?"

# --- Column H width widened to fit the new verification text ---
$ws.Columns("H").ColumnWidth = 16.25

# --- View state: scroll down while keeping the header row frozen ---
$ws.Range("B12").Select()
